# NestedLoopDocument.xlsx — migrate-to-XSSF-backend commit
#
# The underlying OOXML diff this commit produced is dominated by artifacts of
# swapping the export backend (every cellXfs index shifted by one because the
# new writer always inserts its own default style slot at index 0, a
# xl/theme/theme1.xml part appears because the new writer always emits a
# theme, relationship ids were renumbered, attribute ordering changed, etc.).
# None of that is something a user performs through the Excel object model —
# there is no "renumber every style id" or "invent a theme part" action in
# the COM API, and attempting to fake it would not correspond to any real
# spreadsheet edit. The actual, intentional content changes made in the
# workbook are:
#
#   1. A typo fix in the shared string used by cell A48:
#        "gedankeloses" -> "gedankenloses"
#   2. The selection was left on the whole of row 10 (active cell A10)
#      instead of the single cell B39.
#   3. The header/footer font style name changed from "Normal" to "Regular"
#      ("Times New Roman,Normal" -> "Times New Roman,Regular").
#
# These are reproduced below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the typo in the German sentence stored in A48.
$ws.Range("A48").Value = "Das Denken der Gedanken ist ein gedankenloses Denken"

# 2. Leave the selection on row 10 (whole row), active cell A10 — mirrors
#    <selection activeCell="A10" sqref="10:10"/> in the saved worksheet.
$ws.Range("A10").EntireRow.Select()

# 3. Header / footer: font style name "Normal" -> "Regular".
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Página &P'
